$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move content from A3 to B3: copy value, clear old cell, select new cell
$ws.Range("A3").Cut($ws.Range("B3"))

$ws.Range("B3").Select()
